# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.707.97"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "1.678.64"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.68"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3965"
$ws.Range("E8").Value = "  -2.32%  "
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.73"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.410"
$ws.Range("E11").Value = "  -5.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08631"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.28"
$ws.Range("E13").Value = "  -4.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.346"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.807"
$ws.Range("E15").Value = "  -4.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001321"
$ws.Range("E16").Value = "  -2.77%  "
$ws.Range("D17").Value = "1.663.51"
$ws.Range("E17").Value = "  -4.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.74"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07096"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.21"
$ws.Range("E20").Value = "  -4.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.116"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.03"
$ws.Range("E23").Value = "  -2.28%  "
$ws.Range("D24").Value = "24.719.25"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.354"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "23.65"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.778"
$ws.Range("E27").Value = "  -7.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.80"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "150.91"
$ws.Range("E29").Value = "  +3.87%  "
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.749"
$ws.Range("E30").Value = "  -6.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.857"
$ws.Range("E31").Value = "  -7.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.389"
$ws.Range("E32").Value = "  +6.26%  "
$ws.Range("D33").Value = "1.844.12"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08459"
$ws.Range("E34").Value = "  -4.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03090"
$ws.Range("E35").Value = "  -3.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.010"
$ws.Range("E36").Value = "  -3.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.965"
$ws.Range("E37").Value = "  -3.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2804"
$ws.Range("E38").Value = "  -2.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09493"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.50"
$ws.Range("E40").Value = "  -4.11%  "
$ws.Range("E41").Value = "  -4.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.479"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.64"
$ws.Range("E43").Value = "  -3.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.71"
$ws.Range("E44").Value = "  -4.46%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7154"
$ws.Range("E45").Value = "  -3.32%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.579"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.171"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08683"
$ws.Range("E48").Value = "  +4.31%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.339"
$ws.Range("E50").Value = "  -3.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.55"
$ws.Range("E51").Value = "  -2.50%  "
